$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 27-29: mark Status as "Pass" (copy the formatting used by the existing
# "Pass" / date cells in row 26 so the fill/font/border match column-for-column),
# then set the text that actually differs (dates -> "Oct 13th").
$ws.Range("G26:I26").Copy($ws.Range("G27:I27"))
$ws.Range("G26:I26").Copy($ws.Range("G28:I28"))
$ws.Range("G26:I26").Copy($ws.Range("G29:I29"))

$ws.Range("G27").Value2 = "Pass"
$ws.Range("H27").Value2 = "Oct 13th"
$ws.Range("I27").Value2 = "Oct 13th"

$ws.Range("G28").Value2 = "Pass"
$ws.Range("H28").Value2 = "Oct 13th"
$ws.Range("I28").Value2 = "Oct 13th"

$ws.Range("G29").Value2 = "Pass"
$ws.Range("H29").Value2 = "Oct 13th"
$ws.Range("I29").Value2 = "Oct 13th"

# Update the view: scrolled down one row further, and the active selection
# moved from J26 to G28:I28.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
$ws.Range("G28:I28").Select()
